$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents (formulas/values) of the sensitivity analysis table (rows 15-18)
# while preserving existing cell formatting/styles.
$ws.Range("G15:G18").ClearContents()
$ws.Range("H15:H18").ClearContents()
$ws.Range("I15:I18").ClearContents()
$ws.Range("J15:J18").ClearContents()
$ws.Range("K15:K18").ClearContents()
$ws.Range("L15:L18").ClearContents()

# Update the active selection to G14
$ws.Range("G14").Select()
